# Update "合肥-漫展信息" workbook with latest scraped numbers.
# Sheet 1 ("展览") and Sheet 4 ("全部类型") contain the same data table and both
# need to be refreshed identically.

$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Row 2: ticket is no longer sellable -> text replaces the numeric price.
    $ws.Range("G2").Value = "不可售"

    # Row 3: updated "want to go" headcount.
    $ws.Range("F3").Value = 7154

    # Row 4: updated "want to go" headcount.
    $ws.Range("F4").Value = 5172

    # Row 5: updated "want to go" headcount.
    $ws.Range("F5").Value = 76

    # Row 11: updated "want to go" headcount.
    $ws.Range("F11").Value = 88

    # Row 13: updated "want to go" headcount.
    $ws.Range("F13").Value = 630

    # Row 14: updated "want to go" headcount.
    $ws.Range("F14").Value = 199

    # Row 15: updated "want to go" headcount.
    $ws.Range("F15").Value = 51
}
